# Add a new, hidden "Sheet3" containing a single note cell, so the test
# fixture can verify that text extraction skips hidden sheets.

$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore the
# selection afterwards (adding a sheet shifts the active tab).
$previouslyActiveSheetName = $wb.ActiveSheet.Name

# Insert the new sheet after the last existing sheet (so it lands at the
# end of the tab strip, i.e. after Sheet1/Sheet2) instead of the default
# "before the active sheet" placement.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet3"

$ws.Range("A1").Value = "This sheet is hidden."

# Hide the sheet itself.
$ws.Visible = $xlSheetHidden

# Restore the original active sheet/tab selection.
$wb.Worksheets.Item($previouslyActiveSheetName).Activate()
